$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the three new test-plan rows (36-38) describing the new "Edit Quiz" button ---
$ws.Range("A36").Value() = 33
$ws.Range("B36").Value() = "Quizzes Index"
$ws.Range("C36").Value() = "Restricted user is authenticated"
$ws.Range("D36").Value() = '"Edit Quiz" button is not visible'
$ws.Range("E36").Value() = "Same as expected"
$ws.Range("F36").Value() = "Pass"

$ws.Range("A37").Value() = 34
$ws.Range("B37").Value() = "Quizzes Index"
$ws.Range("C37").Value() = "View user is authenticated"
$ws.Range("D37").Value() = '"Edit Quiz" button is not visible'
$ws.Range("E37").Value() = "Same as expected"
$ws.Range("F37").Value() = "Pass"

$ws.Range("A38").Value() = 35
$ws.Range("B38").Value() = "Quizzes Index"
$ws.Range("C38").Value() = "Edit user is authenticated"
$ws.Range("D38").Value() = '"Edit Quiz" button is visible'
$ws.Range("E38").Value() = "Same as expected"
$ws.Range("F38").Value() = "Pass"

# --- Conditional formatting: the pre-existing full-column rules (backed by the
# original 4 dxfs) now only need to cover F36:F38, the rows that were inserted
# in the middle of the old F4:F1048576 range. ---
$allFc = $ws.Cells.FormatConditions
$originalRule = $allFc.Item(1)
$originalRule.ModifyAppliesToRange($ws.Range("F36:F38"))

# The rest of the column (above and below the new rows) needs the same four
# rules re-created (Excel duplicates the dxf/rule set when a CF range is split
# by a row insertion).
$rngRest = $ws.Range("F4:F35")

$fcSuspended = $rngRest.FormatConditions.Add(8, 3, '="Suspended"')
$fcSuspended.Font.Color = 22428
$fcSuspended.Interior.Color = 10284031

$fcNotExecuted = $rngRest.FormatConditions.Add(8, 3, '="Not Executed"')
$fcNotExecuted.Font.Color = 22428
$fcNotExecuted.Interior.Color = 10284031

$fcFail = $rngRest.FormatConditions.Add(8, 3, '="Fail"')
$fcFail.Font.Color = 393372
$fcFail.Interior.Color = 13551615

$fcPass = $rngRest.FormatConditions.Add(8, 3, '="Pass"')
$fcPass.Font.Color = 24832
$fcPass.Interior.Color = 13561798

$rngTail = $ws.Range("F39:F1048576")
$fcSuspended.ModifyAppliesToRange($excel.Union($rngRest, $rngTail))
$fcNotExecuted.ModifyAppliesToRange($excel.Union($rngRest, $rngTail))
$fcFail.ModifyAppliesToRange($excel.Union($rngRest, $rngTail))
$fcPass.ModifyAppliesToRange($excel.Union($rngRest, $rngTail))

# --- Keep the selection in sync with where Excel would have left the cursor
# after typing the last new row. ---
$ws.Range("D38").Select()
